$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell for the Category column (F1). Typed first so "Category"
# becomes the first newly-added shared string.
$ws.Range("F1").Value = "Category"

# Bathroom products: rows 11-19 (towel sets). Filled before Bedroom/Beach so
# "Bathroom" becomes the next new shared string.
for ($r = 11; $r -le 19; $r++) {
    $ws.Cells.Item($r, 6).Value = "Bathroom"
}

# Bedroom products: rows 20-65 (blankets/quilts/sheets).
for ($r = 20; $r -le 65; $r++) {
    $ws.Cells.Item($r, 6).Value = "Bedroom"
}

# Beach products: rows 2-10 (beach towel sets). Filled last so "Beach"
# becomes the final newly-added shared string.
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 6).Value = "Beach"
}

# Give the new column a sensible custom width (closest achievable match).
$ws.Columns.Item(6).ColumnWidth = 21

# Leave the selection on the last data row, matching the saved view state.
$ws.Range("E65").Select() | Out-Null
